$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that contain text which Excel would otherwise auto-convert to
# numbers or dates (I = numeric-looking counts, Y/AA = ISO-style dates).
# Force them to Text format first so the values are written as strings,
# matching the inlineStr/string cell type in the target workbook.
$ws.Range("I3:I6").NumberFormat = "@"
$ws.Range("Y3:Y6").NumberFormat = "@"
$ws.Range("AA3:AA6").NumberFormat = "@"

# ---------------- Row 3 ----------------
$ws.Range("A3").Value = 131106007
$ws.Range("B3").Value = 79243
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("I3").Value = "100"
$ws.Range("J3").Value = "bålar"
$ws.Range("P3").Value = "Svartåsen, Mpd"
$ws.Range("Q3").Value = 612036
$ws.Range("R3").Value = 6945843
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Västernorrland"
$ws.Range("U3").Value = "Timrå"
$ws.Range("V3").Value = "Medelpad"
$ws.Range("W3").Value = "Ljustorp"
$ws.Range("X3").Value = "2025_1183"
$ws.Range("Y3").Value = "2025-09-16"
$ws.Range("Z3").Value = "08:56"
$ws.Range("AA3").Value = "2025-09-16"
$ws.Range("AB3").Value = "08:56"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AW3").Value = "David Isaksson"
$ws.Range("AX3").Value = "Måns Svensson"
$ws.Range("AY3").Value = "Kustpaketet"

# ---------------- Row 4 ----------------
$ws.Range("A4").Value = 131106010
$ws.Range("B4").Value = 79243
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("P4").Value = "Svartåsen, Mpd"
$ws.Range("Q4").Value = 612376
$ws.Range("R4").Value = 6945396
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Västernorrland"
$ws.Range("U4").Value = "Timrå"
$ws.Range("V4").Value = "Medelpad"
$ws.Range("W4").Value = "Ljustorp"
$ws.Range("X4").Value = "2025_1180"
$ws.Range("Y4").Value = "2025-09-16"
$ws.Range("Z4").Value = "08:22"
$ws.Range("AA4").Value = "2025-09-16"
$ws.Range("AB4").Value = "08:22"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "David Isaksson"
$ws.Range("AX4").Value = "Jennifer Lehikoinen, Måns Svensson"
$ws.Range("AY4").Value = "Kustpaketet"

# ---------------- Row 5 ----------------
$ws.Range("A5").Value = 131106008
$ws.Range("B5").Value = 79243
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("I5").Value = "10"
$ws.Range("J5").Value = "bålar"
$ws.Range("P5").Value = "Svartåsen, Mpd"
$ws.Range("Q5").Value = 612048
$ws.Range("R5").Value = 6945825
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Västernorrland"
$ws.Range("U5").Value = "Timrå"
$ws.Range("V5").Value = "Medelpad"
$ws.Range("W5").Value = "Ljustorp"
$ws.Range("X5").Value = "2025_1182"
$ws.Range("Y5").Value = "2025-09-16"
$ws.Range("Z5").Value = "08:51"
$ws.Range("AA5").Value = "2025-09-16"
$ws.Range("AB5").Value = "08:51"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = "David Isaksson"
$ws.Range("AX5").Value = "Jennifer Lehikoinen, Måns Svensson"
$ws.Range("AY5").Value = "Kustpaketet"

# ---------------- Row 6 ----------------
$ws.Range("A6").Value = 131106009
$ws.Range("B6").Value = 79862
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6453
$ws.Range("F6").Value = "Vedskivlav"
$ws.Range("G6").Value = "Hertelidea botryosa"
$ws.Range("H6").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("I6").Value = "8"
$ws.Range("J6").Value = "cm²"
$ws.Range("P6").Value = "Svartåsen, Mpd"
$ws.Range("Q6").Value = 612057
$ws.Range("R6").Value = 6945797
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Västernorrland"
$ws.Range("U6").Value = "Timrå"
$ws.Range("V6").Value = "Medelpad"
$ws.Range("W6").Value = "Ljustorp"
$ws.Range("X6").Value = "2025_1181"
$ws.Range("Y6").Value = "2025-09-16"
$ws.Range("Z6").Value = "08:47"
$ws.Range("AA6").Value = "2025-09-16"
$ws.Range("AB6").Value = "08:47"
$ws.Range("AC6").Value = "tallstubbe"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "David Isaksson"
$ws.Range("AX6").Value = "Måns Svensson"
$ws.Range("AY6").Value = "Kustpaketet"
